# AddictO Organisation Defs - add "Industry" branch (E-cigarette industry,
# Independent e-cigarette industry, Independent e-cigarette company) to the
# Classes sheet, inserted right after the existing "Organisation" header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classes")

# ---------------------------------------------------------------------
# 1. Insert 4 new blank rows starting at row 10 (this pushes the existing
#    "Alcohol Industy" row, and everything below it, down by 4 rows).
# ---------------------------------------------------------------------
$ws.Range("A10:S13").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2. Row 10: new top-level "Industry" category header (mirrors the
#    existing "Organisation" header that sits in row 9).
# ---------------------------------------------------------------------
$ws.Cells.Item(10, 2).Value = "Industry"
$ws.Cells.Item(10, 3).Value = "An aggregate of organisations that manufacture, prepare for sale, distribute, market or promote products or services."

# ---------------------------------------------------------------------
# 3. Row 11: "E-cigarette industry"
# ---------------------------------------------------------------------
$ws.Cells.Item(11, 2).Value  = "E-cigarette industry"
$ws.Cells.Item(11, 3).Value  = "An industry in which the product is electronic cigarettes or e-liquid."
$ws.Cells.Item(11, 4).Value  = "Industry "
$ws.Cells.Item(11, 6).Value  = "Product"
$ws.Cells.Item(11, 8).Value  = "I have added here promotion of - this sets industry apart from manufactureres of but also includes these types of company. "
$ws.Cells.Item(11, 9).Value  = "E-Cig industry"
$ws.Cells.Item(11, 14).Value = 1
$ws.Cells.Item(11, 15).Value = "SC"
$ws.Cells.Item(11, 16).Value = "Proposed"

# ---------------------------------------------------------------------
# 4. Row 12: "Independent e-cigarette industry"
# ---------------------------------------------------------------------
$ws.Cells.Item(12, 2).Value = "Independent e-cigarette industry"
$ws.Cells.Item(12, 3).Value = "An e-cigarette industry that has only independent e-cigarette companies as a part."

# ---------------------------------------------------------------------
# 5. Row 13: "Independent e-cigarette company"
# ---------------------------------------------------------------------
$ws.Cells.Item(13, 2).Value  = "Independent e-cigarette company"
$ws.Cells.Item(13, 3).Value  = "An e-cigarette company that has no commercial or financial association with a tobacco company. "
$ws.Cells.Item(13, 6).Value  = "Product"
$ws.Cells.Item(13, 10).Value = "(Note: This is intended to include ownership, investment or sharing of resources.)"
$ws.Cells.Item(13, 14).Value = 1
$ws.Cells.Item(13, 16).Value = "Proposed"

# ---------------------------------------------------------------------
# 6. Formatting: these four new rows follow the same visual treatment as
#    row 9 (the "Organisation" header) - red font, wrapped text, top
#    vertical alignment, left horizontal alignment for the used cells.
# ---------------------------------------------------------------------
$redFont = 255  # pure red (BGR-ordered COM integer == RGB(255,0,0))

$dataRange = $ws.Range("A10:S13")
$dataRange.Font.Color = $redFont
$dataRange.WrapText = $true
$dataRange.VerticalAlignment = -4160   # xlTop

$usedRange = $ws.Range("B10:O13")
$usedRange.HorizontalAlignment = -4131 # xlLeft

# Row heights (matches auto-fit sizing Excel produced for the wrapped text)
$ws.Rows.Item(10).RowHeight = 29
$ws.Rows.Item(11).RowHeight = 43.5
$ws.Rows.Item(12).RowHeight = 29
$ws.Rows.Item(13).RowHeight = 43.5

# Highlight the "Independent e-cigarette company" definition cell with a
# white fill (kept distinct the way the author did in the source file).
$ws.Cells.Item(13, 3).Interior.Color = 16777215  # white

# Comment cell text doesn't need the strict left-alignment of the other
# columns (kept as vertical-top + wrap only, like the source).
$ws.Cells.Item(13, 10).HorizontalAlignment = -4142 # xlGeneral
